$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'C0603C104M3RACTU'
$ws.Range("B2").Value = 'CAP CER 0.1UF 25V X7R 0603'
$ws.Range("C2").Value = 'C1, C2, C8, C9, C10, C11, C12, C13, C14, C15'
$ws.Range("D2").Value = 'FP-C0603C-CF-MFG'
$ws.Range("E2").Value = 'CMP-2006-03065-2'
$ws.Range("F2").Value = 10
$ws.Range("A3").Value = 'CL10B475KQ8NQNC'
$ws.Range("B3").Value = 'CAP CER 4.7UF 6.3V X7R 0603'
$ws.Range("C3").Value = 'C3, C7, C16'
$ws.Range("D3").Value = 'FP-CL10-IPC_C'
$ws.Range("E3").Value = 'CMP-2000-06096-2'
$ws.Range("F3").Value = 3
$ws.Range("A4").Value = 'CL31A106KBHNNNE'
$ws.Range("B4").Value = 'CAP CER 10UF 50V X5R 1206'
$ws.Range("C4").Value = 'C4'
$ws.Range("D4").Value = 'FP-CL31-IPC_C'
$ws.Range("E4").Value = 'CMP-13271-000016-1'
$ws.Range("F4").Value = 1
$ws.Range("A5").Value = 'GRM1885C1H220JA01D'
$ws.Range("B5").Value = 'Chip Capacitor, 22 pF, +/- 5%, 50 V, -55 to 125 degC, 0603 (1608 Metric), RoHS, Tape and Reel'
$ws.Range("C5").Value = 'C5, C6'
$ws.Range("D5").Value = 'CAPC1608X90X35ML10T15'
$ws.Range("E5").Value = 'CMP-2000-04945-1'
$ws.Range("F5").Value = 2
$ws.Range("A6").Value = '1812L050/30PR'
$ws.Range("B6").Value = 'PTC Resettable Fuse, 0.8 W, 30 V, -40 to 85 degC, 2-Pin SMD, RoHS, Tape and Reel'
$ws.Range("C6").Value = 'F1'
$ws.Range("D6").Value = 'LFUS-1812L05030PR_V'
$ws.Range("E6").Value = 'CMP-2000-05508-1'
$ws.Range("F6").Value = 1
$ws.Range("A7").Value = 'MX34020SF1'
$ws.Range("B7").Value = 'Automotive Connectors Sckt HOUSING 20P'
$ws.Range("C7").Value = 'J1'
$ws.Range("D7").Value = 'MX34020SF1'
$ws.Range("E7").Value = 'MX34020SF1'
$ws.Range("F7").Value = 1
$ws.Range("A8").Value = 'MX34020NF1'
$ws.Range("B8").Value = 'JAE MX34 Series, 2.2mm Pitch 20 Way 2 Row Right Angle PCB Header, Solder Termination, 3A'
$ws.Range("C8").Value = 'J2'
$ws.Range("D8").Value = 'MX34020NF1'
$ws.Range("E8").Value = 'MX34020NF1'
$ws.Range("F8").Value = 1
$ws.Range("A9").Value = '150080SS75000'
$ws.Range("B9").Value = 'SMD mono-color Chip LED, WL-SMCW, Super Red'
$ws.Range("C9").Value = 'LED1'
$ws.Range("D9").Value = '0805_A'
$ws.Range("E9").Value = 'CMP-1426-00010-1'
$ws.Range("F9").Value = 1
$ws.Range("A10").Value = '150080BS75000'
$ws.Range("B10").Value = 'SMD mono-color Chip LED, WL-SMCW, Blue'
$ws.Range("C10").Value = 'LED2'
$ws.Range("D10").Value = '0805_A'
$ws.Range("E10").Value = 'CMP-1426-00008-1'
$ws.Range("F10").Value = 1
$ws.Range("A11").Value = '62201421121'
$ws.Range("B11").Value = 'THT Vertical Pin Header WR-PHD, Pitch 1.27 mm, Dual Row, 14 pins'
$ws.Range("C11").Value = 'P1'
$ws.Range("D11").Value = '62201421121'
$ws.Range("E11").Value = 'CMP-1502-00949-1'
$ws.Range("F11").Value = 1
$ws.Range("A12").Value = '4K7'
$ws.Range("B12").Value = 'Res Thick Film 0603 4.7K Ohm 1% 1/10W 100ppm/C Molded SMD SMD Paper T/R'
$ws.Range("C12").Value = 'R1, R2, R3, R4, R5, R6, R7, R8, R11, R12, R22, R23'
$ws.Range("D12").Value = 'R0603'
$ws.Range("E12").Value = 'RMCF0603FT4K70'
$ws.Range("F12").Value = 12
$ws.Range("A13").Value = 'RC0603FR-07560RL'
$ws.Range("B13").Value = 'Chip Resistor, 560 Ohm, +/- 1%, 0.1 W, -55 to 155 degC, 0603 (1608 Metric), RoHS, Tape and Reel'
$ws.Range("C13").Value = 'R9, R10'
$ws.Range("D13").Value = 'RESC1608X55X25ML10T15'
$ws.Range("E13").Value = 'CMP-1659-00038-1'
$ws.Range("F13").Value = 2
$ws.Range("A14").Value = '76STC02T'
$ws.Range("B14").Value = "SWITCH TOGGLE DIP SPDT 150MA 30V`r`nWIRED:`r`n1  2  3`r`n4  5  6"
$ws.Range("C14").Value = 'SW1'
$ws.Range("D14").Value = '76STC02T'
$ws.Range("E14").Value = '76STC02T'
$ws.Range("F14").Value = 1
$ws.Range("A15").Value = 'UA7805CKTTR'
$ws.Range("B15").Value = "Fixed Positive Voltage Regulator, 7 to 25 V, 1.5 A, 0 to 125 degC, 3-Pin DDPAK (KTT), Green (RoHS`r`n& no Sb/Br), Tape and Reel"
$ws.Range("C15").Value = 'U1'
$ws.Range("D15").Value = 'KTT0003A_V'
$ws.Range("E15").Value = 'CMP-1685-00014-1'
$ws.Range("F15").Value = 1
$ws.Range("A16").Value = 'MCP2562FD-E/SN'
$ws.Range("B16").Value = '8 SOIC 3.90mm(.150in) TUBECAN Flexible Data Rate Transceiver'
$ws.Range("C16").Value = 'U2'
$ws.Range("D16").Value = 'MCP2561'
$ws.Range("E16").Value = 'MCP2562FD-E/SN'
$ws.Range("F16").Value = 1
$ws.Range("A17").Value = 'MIC5504-3.3YM5-TR'
$ws.Range("B17").Value = 'IC REG LINEAR 3.3V 300MA SOT23-5'
$ws.Range("C17").Value = 'U3'
$ws.Range("D17").Value = 'FP-SOT23-5LD-PL-1-MFG'
$ws.Range("E17").Value = 'CMP-2000-07604-2'
$ws.Range("F17").Value = 1
$ws.Range("A18").Value = 'LM339DR2G'
$ws.Range("B18").Value = 'Single Supply Quad Comparators, 0 to 70 degC, 14-Pin SOIC, Pb-Free, Tape and Reel'
$ws.Range("C18").Value = 'U4'
$ws.Range("D18").Value = 'ONSC-SOIC-14-751A-03_V'
$ws.Range("E18").Value = 'CMP-1305-00036-1'
$ws.Range("F18").Value = 1
$ws.Range("A19").Value = 'STM32F103T6U6A'
$ws.Range("B19").Value = 'ARM Cortex-M3 32-bit MCU, 32 KB Flash, 10 KB Internal RAM, 26 I/Os, 36-pin VFQFPN, -40 to 85 degC, Tray'
$ws.Range("C19").Value = 'U5'
$ws.Range("D19").Value = 'STM-VFQFPN36_N'
$ws.Range("E19").Value = 'CMP-0237-00043-3'
$ws.Range("F19").Value = 1
$ws.Range("A20").Value = 'ABM7-8.000MHZ-D2Y-T'
$ws.Range("B20").Value = 'Microprocessor Crystal, 8 MHz, 18 PF, -40 to 85 degC, 2-Pin SMD, RoHS, Tape and Reel'
$ws.Range("C20").Value = 'X1'
$ws.Range("D20").Value = 'ABRA-ABM7-2_V'
$ws.Range("E20").Value = 'CMP-2000-05034-1'
$ws.Range("F20").Value = 1

$ws.PageSetup.Zoom = 49
